# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-57, columns E:G) gets re-sorted from
# descending chronological order (2003 .. 1610) to ascending order
# (1610 .. 2003). Column G is updated uniformly to 781242, while column F
# stays tied to its period value (it just lands on a different row now).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order (was descending 2003 -> 1610).
$periods = @("1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003")

# Column F value associated with each period - unchanged per period, only
# the row it lands on moves because of the re-sort.
$fvals = @(27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249)

# Column G - now a single uniform value for every period.
$gval = 781242

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $fvals[$i]
    $ws.Cells.Item($row, 7).Value = $gval
}
